$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new weekly snapshot row (row 9) mirroring the style of the
# prior rows: A = snapshot_date, B = median_event_date (both dates).
$ws.Range("A9").Value = 46033
$ws.Range("B9").Value = 48994

# Match the date formatting used by the existing rows (same style as A8:B8).
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122) # xlPasteFormats
